$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185:284 down to 186:285
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record's data.
# Columns A,B,C,E,F,G,N,Q,R are constant across this entire dataset.
$ws.Cells.Item(185, 1).Value = 10
$ws.Cells.Item(185, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(185, 3).Value = "La Araucanía"
$ws.Cells.Item(185, 4).Value = 44830
$ws.Cells.Item(185, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 100112043
$ws.Cells.Item(185, 7).Value = "Pepino dulce"
$ws.Cells.Item(185, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 350
$ws.Cells.Item(185, 11).Value = 18000
$ws.Cells.Item(185, 12).Value = 19000
$ws.Cells.Item(185, 13).Value = 18571
$ws.Cells.Item(185, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(185, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(185, 16).Value = 1032
$ws.Cells.Item(185, 17).Value = 18
$ws.Cells.Item(185, 18).Value = "Hortaliza"
